# EI Variable Installments T2 scenarios
#
# Adds a new "waittopageload1" / 2000 step row to the "Edit Repayment
# Schedule" sheet (inserted above the existing "clickonsubmit" row, i.e.
# new row 6; everything below shifts down by one), and leaves that sheet
# as the active sheet/tab with A6:B6 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit Repayment Schedule")

# Insert a new blank row at row 6, pushing the existing rows 6-12 down to 7-13.
$ws.Rows.Item(6).Insert(-4121)  # xlShiftDown

# The newly inserted row already inherits column A's style from the row
# above (matches the other label cells). Column B needs the numeric
# "2000"-style formatting used elsewhere (e.g. row 3), so copy that over.
$ws.Range("B3").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = "waittopageload1"
$ws.Range("B6").Value = 2000

# Make "Edit Repayment Schedule" the active sheet/tab with A6:B6 selected.
$ws.Activate()
$ws.Range("A6:B6").Select()

$wb.Save()
